$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the formatting of row 18 (a "normal" data row) down onto the
#     six rows that follow the current last row (21) so the new rows get
#     the same style indices used throughout the table.
$ws.Range("A18:S18").Copy()
$ws.Range("A21:S26").PasteSpecial(-4122)
for ($r = 21; $r -le 26; $r++) {
  $ws.Rows.Item($r).RowHeight = 27
}

# --- Row 21: "padding"
$ws.Range("A21").Value = "padding"
$ws.Range("B21").Value = 0.15
$ws.Range("C21").Value = 0.45
$ws.Range("D21").Value = 0.5
$ws.Range("E21").Formula = "=IF((`$E`$1/100*B21/16)+C21<D21,D21,IF((`$E`$1/100*B21/16)+C21>J21,J21,(`$E`$1/100*B21/16)+C21))"
$ws.Range("F21").Formula = "=IF((`$F`$1/100*B21/16)+C21<D21,D21,IF((`$F`$1/100*B21/16)+C21>J21,J21,(`$F`$1/100*B21/16)+C21))"
$ws.Range("G21").Formula = "=IF((`$G`$1/100*B21/16)+C21<D21,D21,IF((`$G`$1/100*B21/16)+C21>J21,J21,(`$G`$1/100*B21/16)+C21))"
$ws.Range("H21").Formula = "=IF((`$H`$1/100*B21/16)+C21<D21,D21,IF((`$H`$1/100*B21/16)+C21>J21,J21,(`$H`$1/100*B21/16)+C21))"
$ws.Range("I21").Formula = "=IF((`$I`$1/100*B21/16)+C21<D21,D21,IF((`$I`$1/100*B21/16)+C21>J21,J21,(`$I`$1/100*B21/16)+C21))"
$ws.Range("J21").Value = 0.63
$ws.Range("K21").Value = 0.8
$ws.Range("L21").Formula = "=E21*16"
$ws.Range("M21").Formula = "=F21*16"
$ws.Range("N21").Formula = "=G21*16"
$ws.Range("O21").Formula = "=H21*16"
$ws.Range("P21").Formula = "=I21*16"
$ws.Range("S21").Formula = '="font-size: clamp("&D21&"rem, "&B21&"vw + "&C21&"rem, "&J21&"rem);"'

# --- Row 22: "width"
$ws.Range("A22").Value = "width"
$ws.Range("B22").Value = 1.1
$ws.Range("C22").Value = 2.67
$ws.Range("D22").Value = 3
$ws.Range("E22").Formula = "=IF((`$E`$1/100*B22/16)+C22<D22,D22,IF((`$E`$1/100*B22/16)+C22>J22,J22,(`$E`$1/100*B22/16)+C22))"
$ws.Range("F22").Formula = "=IF((`$F`$1/100*B22/16)+C22<D22,D22,IF((`$F`$1/100*B22/16)+C22>J22,J22,(`$F`$1/100*B22/16)+C22))"
$ws.Range("G22").Formula = "=IF((`$G`$1/100*B22/16)+C22<D22,D22,IF((`$G`$1/100*B22/16)+C22>J22,J22,(`$G`$1/100*B22/16)+C22))"
$ws.Range("H22").Formula = "=IF((`$H`$1/100*B22/16)+C22<D22,D22,IF((`$H`$1/100*B22/16)+C22>J22,J22,(`$H`$1/100*B22/16)+C22))"
$ws.Range("I22").Formula = "=IF((`$I`$1/100*B22/16)+C22<D22,D22,IF((`$I`$1/100*B22/16)+C22>J22,J22,(`$I`$1/100*B22/16)+C22))"
$ws.Range("J22").Value = 4
$ws.Range("K22").Value = 2.5
$ws.Range("L22").Formula = "=E22*16"
$ws.Range("M22").Formula = "=F22*16"
$ws.Range("N22").Formula = "=G22*16"
$ws.Range("O22").Formula = "=H22*16"
$ws.Range("P22").Formula = "=I22*16"
$ws.Range("S22").Formula = '="font-size: clamp("&D22&"rem, "&B22&"vw + "&C22&"rem, "&J22&"rem);"'

# --- Column widths (character widths map to XML width = chars + 5/6)
$ws.Columns.Item(1).ColumnWidth = 7.9166666666666667
$ws.Range("B1:Q1").EntireColumn.ColumnWidth = 7.9166666666666667
$ws.Columns.Item(18).ColumnWidth = 8.2916666666666667
$ws.Columns.Item(19).ColumnWidth = 45.4166666666666667

# --- Selection / view bookkeeping (cosmetic, mirrors the recorded edit)
$ws.Range("S22").Select()
